$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 614 ("「コーヒーは急には飲まれない」..." entry).
# This shifts all subsequent rows up by one, matching the target diff
# (row 615 becomes 614, ..., row 773 becomes 772).
$ws.Rows.Item(614).Delete()
